$wb = $excel.ActiveWorkbook

# --- Variables sheet: change "country" row's valueType from "text" to "integer" ---
$wsVariables = $wb.Worksheets.Item("Variables")
$wsVariables.Range("B4").Value = "integer"

# --- Categories sheet: change "name" column from text codes to sequential integers ---
$wsCategories = $wb.Worksheets.Item("Categories")
$wsCategories.Range("B2").Value = 0
$wsCategories.Range("B3").Value = 1
$wsCategories.Range("B4").Value = 2
$wsCategories.Range("B5").Value = 3

# --- Selection / active sheet state ---
$wsVariables.Range("D15").Select()
$wsCategories.Select()
$wsCategories.Range("B6").Select()
